{"js": "// Replace the date line and each three-digit-by-one-digit multiplication\n// answer in the table with the new values described by the commit diff.\n// Every \"old\" string below is unique within the document, so a simple\n// matchCase/wholeMatch-free search-and-replace on each pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"2025-06-14 Saturday\", \"2025-06-15 Sunday\"],\n  [\"835\u00d75=4175\", \"165\u00d78=1320\"],\n  [\"614\u00d74=2456\", \"331\u00d76=1986\"],\n  [\"931\u00d73=2793\", \"823\u00d75=4115\"],\n  [\"556\u00d74=2224\", \"913\u00d76=5478\"],\n  [\"350\u00d73=1050\", \"971\u00d78=7768\"],\n  [\"966\u00d73=2898\", \"834\u00d73=2502\"],\n  [\"396\u00d75=1980\", \"603\u00d78=4824\"],\n  [\"136\u00d74=544\", \"583\u00d79=5247\"],\n  [\"754\u00d73=2262\", \"622\u00d76=3732\"],\n  [\"217\u00d79=1953\", \"649\u00d74=2596\"],\n  [\"257\u00d75=1285\", \"586\u00d72=1172\"],\n  [\"233\u00d77=1631\", \"789\u00d75=3945\"],\n  [\"265\u00d76=1590\", \"477\u00d79=4293\"],\n  [\"690\u00d77=4830\", \"453\u00d78=3624\"],\n  [\"384\u00d76=2304\", \"677\u00d78=5416\"],\n  [\"160\u00d79=1440\", \"236\u00d73=708\"],\n  [\"733\u00d72=1466\", \"657\u00d78=5256\"],\n  [\"972\u00d72=1944\", \"678\u00d73=2034\"],\n  [\"330\u00d76=1980\", \"905\u00d74=3620\"],\n  [\"943\u00d74=3772\", \"816\u00d72=1632\"],\n  [\"866\u00d75=4330\", \"468\u00d77=3276\"],\n  [\"152\u00d72=304\", \"708\u00d79=6372\"],\n  [\"153\u00d76=918\", \"587\u00d79=5283\"],\n  [\"883\u00d74=3532\", \"377\u00d73=1131\"],\n  [\"335\u00d77=2345\", \"307\u00d79=2763\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each three-digit-by-one-digit multiplication\n# answer in the table with the new values described by the commit diff.\n# Every \"old\" string is unique within the document, so Find/Replace with\n# wdReplaceAll for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-14 Saturday\", \"2025-06-15 Sunday\"),\n    @(\"835\u00d75=4175\", \"165\u00d78=1320\"),\n    @(\"614\u00d74=2456\", \"331\u00d76=1986\"),\n    @(\"931\u00d73=2793\", \"823\u00d75=4115\"),\n    @(\"556\u00d74=2224\", \"913\u00d76=5478\"),\n    @(\"350\u00d73=1050\", \"971\u00d78=7768\"),\n    @(\"966\u00d73=2898\", \"834\u00d73=2502\"),\n    @(\"396\u00d75=1980\", \"603\u00d78=4824\"),\n    @(\"136\u00d74=544\", \"583\u00d79=5247\"),\n    @(\"754\u00d73=2262\", \"622\u00d76=3732\"),\n    @(\"217\u00d79=1953\", \"649\u00d74=2596\"),\n    @(\"257\u00d75=1285\", \"586\u00d72=1172\"),\n    @(\"233\u00d77=1631\", \"789\u00d75=3945\"),\n    @(\"265\u00d76=1590\", \"477\u00d79=4293\"),\n    @(\"690\u00d77=4830\", \"453\u00d78=3624\"),\n    @(\"384\u00d76=2304\", \"677\u00d78=5416\"),\n    @(\"160\u00d79=1440\", \"236\u00d73=708\"),\n    @(\"733\u00d72=1466\", \"657\u00d78=5256\"),\n    @(\"972\u00d72=1944\", \"678\u00d73=2034\"),\n    @(\"330\u00d76=1980\", \"905\u00d74=3620\"),\n    @(\"943\u00d74=3772\", \"816\u00d72=1632\"),\n    @(\"866\u00d75=4330\", \"468\u00d77=3276\"),\n    @(\"152\u00d72=304\", \"708\u00d79=6372\"),\n    @(\"153\u00d76=918\", \"587\u00d79=5283\"),\n    @(\"883\u00d74=3532\", \"377\u00d73=1131\"),\n    @(\"335\u00d77=2345\", \"307\u00d79=2763\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
